$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose match data was corrected/reordered ---
# Row 58
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "spain"
$ws.Range("C58").Value = "laliga"
$ws.Range("D58").Value = "2023-2024"
$ws.Range("E58").Value = 45193.77083333334
$ws.Range("F58").Value = "Las Palmas"
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = "Granada CF"
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2.14
$ws.Range("K58").Value = "11/09/2023 13:19"
$ws.Range("L58").Value = 1.95
$ws.Range("M58").Value = "24/09/2023 18:19"
$ws.Range("N58").Value = 3.21
$ws.Range("O58").Value = "11/09/2023 13:19"
$ws.Range("P58").Value = 3.73
$ws.Range("Q58").Value = "24/09/2023 18:27"
$ws.Range("R58").Value = 3.61
$ws.Range("S58").Value = "11/09/2023 13:19"
$ws.Range("T58").Value = 4.1
$ws.Range("U58").Value = "24/09/2023 18:27"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/spain/laliga/las-palmas-granada-cf/tWsBDE3N/"

# Row 59
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "spain"
$ws.Range("C59").Value = "laliga"
$ws.Range("D59").Value = "2023-2024"
$ws.Range("E59").Value = 45193.77083333334
$ws.Range("F59").Value = "Betis"
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = "Cadiz CF"
$ws.Range("I59").Value = 1
$ws.Range("J59").Value = 1.87
$ws.Range("K59").Value = "05/09/2023 12:02"
$ws.Range("L59").Value = 1.81
$ws.Range("M59").Value = "24/09/2023 18:26"
$ws.Range("N59").Value = 3.56
$ws.Range("O59").Value = "05/09/2023 12:02"
$ws.Range("P59").Value = 3.66
$ws.Range("Q59").Value = "24/09/2023 18:26"
$ws.Range("R59").Value = 4.54
$ws.Range("S59").Value = "05/09/2023 12:02"
$ws.Range("T59").Value = 5.07
$ws.Range("U59").Value = "24/09/2023 18:26"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/spain/laliga/betis-cadiz/IicoJIZo/"

# Row 63
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = "spain"
$ws.Range("C63").Value = "laliga"
$ws.Range("D63").Value = "2023-2024"
$ws.Range("E63").Value = 45196.79166666666
$ws.Range("F63").Value = "Ath Bilbao"
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = "Getafe"
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 1.71
$ws.Range("K63").Value = "17/09/2023 09:02"
$ws.Range("L63").Value = 1.53
$ws.Range("M63").Value = "27/09/2023 18:31"
$ws.Range("N63").Value = 3.42
$ws.Range("O63").Value = "17/09/2023 09:02"
$ws.Range("P63").Value = 4.06
$ws.Range("Q63").Value = "27/09/2023 18:49"
$ws.Range("R63").Value = 5.44
$ws.Range("S63").Value = "17/09/2023 09:02"
$ws.Range("T63").Value = 7.73
$ws.Range("U63").Value = "27/09/2023 18:49"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/spain/laliga/ath-bilbao-getafe/zgsFCYIT/"

# Row 64
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = "spain"
$ws.Range("C64").Value = "laliga"
$ws.Range("D64").Value = "2023-2024"
$ws.Range("E64").Value = 45196.79166666666
$ws.Range("F64").Value = "Villarreal"
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = "Girona"
$ws.Range("I64").Value = 2
$ws.Range("J64").Value = 1.71
$ws.Range("K64").Value = "17/09/2023 09:02"
$ws.Range("L64").Value = 2.17
$ws.Range("M64").Value = "27/09/2023 18:51"
$ws.Range("N64").Value = 4.14
$ws.Range("O64").Value = "17/09/2023 09:02"
$ws.Range("P64").Value = 3.72
$ws.Range("Q64").Value = "27/09/2023 18:51"
$ws.Range("R64").Value = 4.72
$ws.Range("S64").Value = "17/09/2023 09:02"
$ws.Range("T64").Value = 3.42
$ws.Range("U64").Value = "27/09/2023 18:51"
$ws.Range("V64").Value = "https://www.betexplorer.com/football/spain/laliga/villarreal-girona/80EuTg3A/"

# Row 65
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = "spain"
$ws.Range("C65").Value = "laliga"
$ws.Range("D65").Value = "2023-2024"
$ws.Range("E65").Value = 45196.79166666666
$ws.Range("F65").Value = "Real Madrid"
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = "Las Palmas"
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1.2
$ws.Range("K65").Value = "23/09/2023 09:28"
$ws.Range("L65").Value = 1.18
$ws.Range("M65").Value = "27/09/2023 18:29"
$ws.Range("N65").Value = 6.76
$ws.Range("O65").Value = "23/09/2023 09:28"
$ws.Range("P65").Value = 8
$ws.Range("Q65").Value = "27/09/2023 18:29"
$ws.Range("R65").Value = 11.3
$ws.Range("S65").Value = "23/09/2023 09:28"
$ws.Range("T65").Value = 16.5
$ws.Range("U65").Value = "27/09/2023 18:29"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/spain/laliga/real-madrid-las-palmas/GQHmRXXM/"

# Row 88
$ws.Range("A88").Value = 87
$ws.Range("B88").Value = "spain"
$ws.Range("C88").Value = "laliga"
$ws.Range("D88").Value = "2023-2024"
$ws.Range("E88").Value = 45207.77083333334
$ws.Range("F88").Value = "Alaves"
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = "Betis"
$ws.Range("I88").Value = 1
$ws.Range("J88").Value = 2.8
$ws.Range("K88").Value = "01/10/2023 20:24"
$ws.Range("L88").Value = 2.57
$ws.Range("M88").Value = "08/10/2023 18:28"
$ws.Range("N88").Value = 3.04
$ws.Range("O88").Value = "01/10/2023 20:24"
$ws.Range("P88").Value = 3.25
$ws.Range("Q88").Value = "08/10/2023 18:28"
$ws.Range("R88").Value = 2.72
$ws.Range("S88").Value = "01/10/2023 20:24"
$ws.Range("T88").Value = 3.06
$ws.Range("U88").Value = "08/10/2023 18:22"
$ws.Range("V88").Value = "https://www.betexplorer.com/football/spain/laliga/alaves-betis/YNPlfW19/"

# Row 89
$ws.Range("A89").Value = 88
$ws.Range("B89").Value = "spain"
$ws.Range("C89").Value = "laliga"
$ws.Range("D89").Value = "2023-2024"
$ws.Range("E89").Value = 45207.77083333334
$ws.Range("F89").Value = "Celta Vigo"
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = "Getafe"
$ws.Range("I89").Value = 2
$ws.Range("J89").Value = 1.92
$ws.Range("K89").Value = "28/09/2023 15:02"
$ws.Range("L89").Value = 2.04
$ws.Range("M89").Value = "08/10/2023 18:29"
$ws.Range("N89").Value = 3.25
$ws.Range("O89").Value = "28/09/2023 15:02"
$ws.Range("P89").Value = 3.29
$ws.Range("Q89").Value = "08/10/2023 18:27"
$ws.Range("R89").Value = 4.85
$ws.Range("S89").Value = "28/09/2023 15:02"
$ws.Range("T89").Value = 4.38
$ws.Range("U89").Value = "08/10/2023 18:29"
$ws.Range("V89").Value = "https://www.betexplorer.com/football/spain/laliga/celta-vigo-getafe/0ARtdhXd/"

# --- Append new rows 101-110 (copy formatting from row 100 template) ---
$ws.Range("A100:V100").Copy()
$ws.Range("A101:V110").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 101
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = "spain"
$ws.Range("C101").Value = "laliga"
$ws.Range("D101").Value = "2023-2024"
$ws.Range("E101").Value = 45226.875
$ws.Range("F101").Value = "Girona"
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = "Celta Vigo"
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 2.01
$ws.Range("K101").Value = "10/10/2023 14:02"
$ws.Range("L101").Value = 1.75
$ws.Range("M101").Value = "27/10/2023 20:59"
$ws.Range("N101").Value = 3.56
$ws.Range("O101").Value = "10/10/2023 14:02"
$ws.Range("P101").Value = 4.12
$ws.Range("Q101").Value = "27/10/2023 20:59"
$ws.Range("R101").Value = 3.89
$ws.Range("S101").Value = "10/10/2023 14:02"
$ws.Range("T101").Value = 4.63
$ws.Range("U101").Value = "27/10/2023 20:59"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/spain/laliga/girona-celta-vigo/UVYk0bUa/"

# Row 102
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = "spain"
$ws.Range("C102").Value = "laliga"
$ws.Range("D102").Value = "2023-2024"
$ws.Range("E102").Value = 45227.58333333334
$ws.Range("F102").Value = "Almeria"
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = "Las Palmas"
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = 2.28
$ws.Range("K102").Value = "10/10/2023 14:29"
$ws.Range("L102").Value = 2.31
$ws.Range("M102").Value = "28/10/2023 13:27"
$ws.Range("N102").Value = 3.35
$ws.Range("O102").Value = "10/10/2023 14:29"
$ws.Range("P102").Value = 3.3
$ws.Range("Q102").Value = "28/10/2023 13:27"
$ws.Range("R102").Value = 3.15
$ws.Range("S102").Value = "10/10/2023 14:29"
$ws.Range("T102").Value = 3.42
$ws.Range("U102").Value = "28/10/2023 13:25"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/spain/laliga/almeria-las-palmas/rZCKpJUC/"

# Row 103
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = "spain"
$ws.Range("C103").Value = "laliga"
$ws.Range("D103").Value = "2023-2024"
$ws.Range("E103").Value = 45227.67708333334
$ws.Range("F103").Value = "Barcelona"
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = "Real Madrid"
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = 1.94
$ws.Range("K103").Value = "10/10/2023 14:02"
$ws.Range("L103").Value = 2.58
$ws.Range("M103").Value = "28/10/2023 16:13"
$ws.Range("N103").Value = 3.89
$ws.Range("O103").Value = "10/10/2023 14:02"
$ws.Range("P103").Value = 3.49
$ws.Range("Q103").Value = "28/10/2023 16:14"
$ws.Range("R103").Value = 3.8
$ws.Range("S103").Value = "10/10/2023 14:02"
$ws.Range("T103").Value = 2.83
$ws.Range("U103").Value = "28/10/2023 16:13"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/spain/laliga/barcelona-real-madrid/OUk78MMJ/"

# Row 104
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = "spain"
$ws.Range("C104").Value = "laliga"
$ws.Range("D104").Value = "2023-2024"
$ws.Range("E104").Value = 45227.77083333334
$ws.Range("F104").Value = "Mallorca"
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = "Getafe"
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2.19
$ws.Range("K104").Value = "10/10/2023 14:02"
$ws.Range("L104").Value = 2.33
$ws.Range("M104").Value = "28/10/2023 18:05"
$ws.Range("N104").Value = 2.95
$ws.Range("O104").Value = "10/10/2023 14:02"
$ws.Range("P104").Value = 2.96
$ws.Range("Q104").Value = "28/10/2023 18:05"
$ws.Range("R104").Value = 4.19
$ws.Range("S104").Value = "10/10/2023 14:02"
$ws.Range("T104").Value = 3.9
$ws.Range("U104").Value = "28/10/2023 18:24"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/spain/laliga/mallorca-getafe/OhYgaIq6/"

# Row 105
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = "spain"
$ws.Range("C105").Value = "laliga"
$ws.Range("D105").Value = "2023-2024"
$ws.Range("E105").Value = 45227.875
$ws.Range("F105").Value = "Cadiz CF"
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = "Sevilla"
$ws.Range("I105").Value = 2
$ws.Range("J105").Value = 2.93
$ws.Range("K105").Value = "10/10/2023 14:02"
$ws.Range("L105").Value = 3.9
$ws.Range("M105").Value = "28/10/2023 20:58"
$ws.Range("N105").Value = 3.18
$ws.Range("O105").Value = "10/10/2023 14:02"
$ws.Range("P105").Value = 3.49
$ws.Range("Q105").Value = "28/10/2023 20:58"
$ws.Range("R105").Value = 2.65
$ws.Range("S105").Value = "10/10/2023 14:02"
$ws.Range("T105").Value = 2.08
$ws.Range("U105").Value = "28/10/2023 20:43"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/spain/laliga/cadiz-sevilla/fDtw30ps/"

# Row 106
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = "spain"
$ws.Range("C106").Value = "laliga"
$ws.Range("D106").Value = "2023-2024"
$ws.Range("E106").Value = 45228.58333333334
$ws.Range("F106").Value = "Betis"
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = "Osasuna"
$ws.Range("I106").Value = 1
$ws.Range("J106").Value = 1.99
$ws.Range("K106").Value = "10/10/2023 14:02"
$ws.Range("L106").Value = 1.99
$ws.Range("M106").Value = "29/10/2023 13:59"
$ws.Range("N106").Value = 3.42
$ws.Range("O106").Value = "10/10/2023 14:02"
$ws.Range("P106").Value = 3.47
$ws.Range("Q106").Value = "29/10/2023 13:59"
$ws.Range("R106").Value = 4.17
$ws.Range("S106").Value = "10/10/2023 14:02"
$ws.Range("T106").Value = 4.21
$ws.Range("U106").Value = "29/10/2023 13:59"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/spain/laliga/betis-osasuna/xYgB7txQ/"

# Row 107
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = "spain"
$ws.Range("C107").Value = "laliga"
$ws.Range("D107").Value = "2023-2024"
$ws.Range("E107").Value = 45228.67708333334
$ws.Range("F107").Value = "Rayo Vallecano"
$ws.Range("G107").Value = 2
$ws.Range("H107").Value = "Real Sociedad"
$ws.Range("I107").Value = 2
$ws.Range("J107").Value = 2.86
$ws.Range("K107").Value = "10/10/2023 14:02"
$ws.Range("L107").Value = 3.38
$ws.Range("M107").Value = "29/10/2023 15:55"
$ws.Range("N107").Value = 3.06
$ws.Range("O107").Value = "10/10/2023 14:02"
$ws.Range("P107").Value = 3.16
$ws.Range("Q107").Value = "29/10/2023 16:13"
$ws.Range("R107").Value = 2.65
$ws.Range("S107").Value = "10/10/2023 14:02"
$ws.Range("T107").Value = 2.42
$ws.Range("U107").Value = "29/10/2023 15:55"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/spain/laliga/rayo-vallecano-real-sociedad/Y3us2Kam/"

# Row 108
$ws.Range("A108").Value = 107
$ws.Range("B108").Value = "spain"
$ws.Range("C108").Value = "laliga"
$ws.Range("D108").Value = "2023-2024"
$ws.Range("E108").Value = 45228.77083333334
$ws.Range("F108").Value = "Ath Bilbao"
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = "Valencia"
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 1.69
$ws.Range("K108").Value = "10/10/2023 14:02"
$ws.Range("L108").Value = 1.65
$ws.Range("M108").Value = "29/10/2023 18:29"
$ws.Range("N108").Value = 3.67
$ws.Range("O108").Value = "10/10/2023 14:02"
$ws.Range("P108").Value = 3.97
$ws.Range("Q108").Value = "29/10/2023 18:29"
$ws.Range("R108").Value = 5.04
$ws.Range("S108").Value = "10/10/2023 14:02"
$ws.Range("T108").Value = 5.83
$ws.Range("U108").Value = "29/10/2023 18:29"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/spain/laliga/ath-bilbao-valencia/MPBOqwpJ/"

# Row 109
$ws.Range("A109").Value = 108
$ws.Range("B109").Value = "spain"
$ws.Range("C109").Value = "laliga"
$ws.Range("D109").Value = "2023-2024"
$ws.Range("E109").Value = 45228.875
$ws.Range("F109").Value = "Atl. Madrid"
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = "Alaves"
$ws.Range("I109").Value = 1
$ws.Range("J109").Value = 1.33
$ws.Range("K109").Value = "10/10/2023 14:31"
$ws.Range("L109").Value = 1.45
$ws.Range("M109").Value = "29/10/2023 20:54"
$ws.Range("N109").Value = 4.96
$ws.Range("O109").Value = "10/10/2023 14:31"
$ws.Range("P109").Value = 4.57
$ws.Range("Q109").Value = "29/10/2023 20:59"
$ws.Range("R109").Value = 9.220000000000001
$ws.Range("S109").Value = "10/10/2023 14:31"
$ws.Range("T109").Value = 8.199999999999999
$ws.Range("U109").Value = "29/10/2023 20:59"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/spain/laliga/atl-madrid-alaves/zTybbxaC/"

# Row 110
$ws.Range("A110").Value = 109
$ws.Range("B110").Value = "spain"
$ws.Range("C110").Value = "laliga"
$ws.Range("D110").Value = "2023-2024"
$ws.Range("E110").Value = 45229.875
$ws.Range("F110").Value = "Granada CF"
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = "Villarreal"
$ws.Range("I110").Value = 3
$ws.Range("J110").Value = 3.27
$ws.Range("K110").Value = "10/10/2023 14:34"
$ws.Range("L110").Value = 2.99
$ws.Range("M110").Value = "30/10/2023 20:59"
$ws.Range("N110").Value = 3.51
$ws.Range("O110").Value = "10/10/2023 14:34"
$ws.Range("P110").Value = 3.64
$ws.Range("Q110").Value = "30/10/2023 20:58"
$ws.Range("R110").Value = 2.15
$ws.Range("S110").Value = "10/10/2023 14:34"
$ws.Range("T110").Value = 2.39
$ws.Range("U110").Value = "30/10/2023 20:59"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/spain/laliga/granada-cf-villarreal/juzn1vFg/"
